$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet / update title to reflect new date 2022-09-23
$ws.Name = "Through 2022-09-23"

# Update the "September (through 09-22)" label to "September (through 09-23)"
$ws.Range("A10").Value = "September (through 09-23)"

# Update September row (row 10) values
$ws.Range("B10").Value = 25
$ws.Range("D10").Value = 57
$ws.Range("E10").Value = 46
$ws.Range("F10").Value = 57
$ws.Range("G10").Value = 87
$ws.Range("H10").Value = 134
$ws.Range("I10").Value = 111

# Update Total row (row 11) values
$ws.Range("B11").Value = 219
$ws.Range("D11").Value = 608
$ws.Range("E11").Value = 536
$ws.Range("F11").Value = 406
$ws.Range("G11").Value = 871
$ws.Range("H11").Value = 1204
$ws.Range("I11").Value = 1246
